# Keymap.pptx — "Added: toggle envelope height (shift + h)"
#
# 1) Nudge three existing label/key shapes that sit next to the new
#    "Toggle Envelope height" key-combo graphic (they get slightly
#    repositioned to make room for it).
# 2) Add the new "⇧" key glyph shape + its two-line caption shape.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

$EMU = 12700.0

# --- 1) reposition existing shapes ------------------------------------

$t297 = Get-ShapeById $s 297
$t297.Left = 4106358 / $EMU
$t297.Top  = 3748929 / $EMU

$t207 = Get-ShapeById $s 207
$t207.Left = 4278995 / $EMU
$t207.Top  = 4043086 / $EMU

$t349 = Get-ShapeById $s 349
$t349.Left = 4408719 / $EMU
$t349.Top  = 4062470 / $EMU

# --- 2) add the new "⇧" key glyph -------------------------------------

$shift = $s.Shapes.AddShape(1, 4454008 / $EMU, 3856194 / $EMU, 248786 / $EMU, 246221 / $EMU)
$shift.Name = "Rechteck 361"

$tf = $shift.TextFrame
$tf.WordWrap = $false
$tf.AutoSize = 1

$tr = $tf.TextRange
$tr.Text = "⇧"
$tr.Font.Size = 10
$tr.Font.Name = "Helvetica Neue LT Std 67 Medium Condensed"
$tr.Font.Color.ObjectThemeColor = 9
$tr.Font.Color.TintAndShade = 0
$tr.Font.Color.Brightness = -0.25

# --- 3) add the new two-line caption ----------------------------------

$cap = $s.Shapes.AddShape(1, 4291791 / $EMU, 3913480 / $EMU, 579005 / $EMU, 220573 / $EMU)
$cap.Name = "Rechteck 362"

$ctf = $cap.TextFrame
$ctf.WordWrap = $false
$ctf.AutoSize = 1

$ctr = $ctf.TextRange
$ctr.Text = " Toggle`rEnvelope height"
$ctr.Font.Size = 5
$ctr.Font.Name = "Helvetica Neue LT Std 67 Medium Condensed"
$ctr.Font.Color.ObjectThemeColor = 9
$ctr.Font.Color.TintAndShade = 0
$ctr.Font.Color.Brightness = -0.25

$ctf.Ruler.Levels.Item(1).ParagraphFormat.Alignment = 3
